$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.004.72'
$ws.Range("E2").Value = '  -3.88%  '
$ws.Range("D3").Value = '1.644.25'
$ws.Range("E3").Value = '  -5.62%  '
$ws.Range("D4").Value = '''0.9974'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '''233.53'
$ws.Range("E5").Value = '  -5.35%  '
$ws.Range("D6").Value = '''0.9996'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '''0.4797'
$ws.Range("E7").Value = '  -5.00%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '''0.2599'
$ws.Range("E8").Value = '  -5.20%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.06146'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").Value = '''0.07079'
$ws.Range("E10").Value = '  -2.30%  '
$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '1.643.65'
$ws.Range("E11").Value = '  -5.66%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '''14.66'
$ws.Range("E12").Value = '  -3.12%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.5996'
$ws.Range("E13").Value = '  -8.05%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.366'
$ws.Range("E14").Value = '  -6.72%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '''73.71'
$ws.Range("E15").Value = '  -4.98%  '
$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D16").Value = '''0.9997'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '''0.9986'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '24.995.15'
$ws.Range("E18").Value = '  -3.95%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '''0.000006605'
$ws.Range("E19").Value = '  -3.80%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '''11.28'
$ws.Range("E20").Value = '  -5.22%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '1.850.78'
$ws.Range("E21").Value = '  -6.15%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''4.362'
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '''8.582'
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '''5.250'
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''134.15'
$ws.Range("E25").Value = '  -1.28%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '''14.90'
$ws.Range("E26").Value = '  -2.39%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '''1.386'
$ws.Range("E27").Value = '  -7.63%  '
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '''103.97'
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '''1.648'
$ws.Range("E29").Value = '  -7.57%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '''3.882'
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.07681'
$ws.Range("E31").Value = '  -6.11%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''3.562'
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("B33").Value = 'Frax'
$ws.Range("C33").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D33").Value = '''0.9987'
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.04294'
$ws.Range("E34").Value = '  -8.06%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.567'
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''0.9293'
$ws.Range("E36").Value = '  -6.89%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.5925'
$ws.Range("E37").Value = '  -2.91%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '''2.576'
$ws.Range("E38").Value = '  -7.45%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '''0.8515'
$ws.Range("E39").Value = '  +11.50%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '''0.9988'
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '''0.01516'
$ws.Range("E41").Value = '  -6.58%  '
$ws.Range("B42").Value = 'PaxosStandard'
$ws.Range("C42").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D42").Value = '''0.9989'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '''98.49'
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("D44").Value = '''1.762'
$ws.Range("E44").Value = '  -8.71%  '
$ws.Range("D45").Value = '''0.3704'
$ws.Range("E45").Value = '  -5.65%  '
$ws.Range("D46").Value = '''4.687'
$ws.Range("E46").Value = '  -6.33%  '
$ws.Range("D47").Value = '''0.1096'
$ws.Range("E47").Value = '  -5.94%  '
$ws.Range("D48").Value = '''6.096'
$ws.Range("E48").Value = '  -3.71%  '
$ws.Range("D49").Value = '''0.05201'
$ws.Range("E49").Value = '  -1.90%  '
$ws.Range("D50").Value = '''29.05'
$ws.Range("E50").Value = '  -5.24%  '
$ws.Range("D51").Value = '''0.9988'
$ws.Range("E51").Value = '  -0.32%  '

# Reset style on apostrophe-entered numeric-looking text cells so Excel
# drops the auto-applied quotePrefix formatting (keeps cells style-less,
# matching the original workbook which has no style on these data cells).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
